# Agregando nuevas funciones automatizadas
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The UiPath automation now stamps a run status ("SUCCESS") into column N
# for every processed contact, and logs an error detail into column O
# whenever a step on that row failed.
$ws.Range("N2:N8").Value = "SUCCESS"
$ws.Range("O3").Value = "Error: The UiElement is no longer valid Fuente: Estado"

# Reflect where the user had scrolled/selected when the workbook was last
# saved: column G visible on the left edge, active selection on R16.
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("R16").Select()

# Restore the window geometry recorded for this session.
$excel.ActiveWindow.Left = 3030
$excel.ActiveWindow.Top = 3030
$excel.ActiveWindow.Width = 15375
$excel.ActiveWindow.Height = 7875
